# Auto-generated Excel COM-interop script applying scheduled market-price refresh
# to Sheets/Masamune_Profits.xlsx (per diff). Updates currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H-N) for the affected leve rows across sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 95: Official Strategy Guide | Gyuki Leather Codex
$ws.Cells.Item(95, 8).Value = 37016.668  # H95: 37032.668 -> 37016.668
$ws.Cells.Item(95, 10).Value = 37016.668  # J95: 37032.668 -> 37016.668
$ws.Cells.Item(95, 12).Value = 37016.668  # L95: 37032.668 -> 37016.668
$ws.Cells.Item(95, 14).Value = -42508.668  # N95: -42524.668 -> -42508.668

# Row 107: Another Man's Ink | Enchanted Truegold Ink
$ws.Cells.Item(107, 8).Value = 419.70587  # H107: 8017.231 -> 419.70587
$ws.Cells.Item(107, 9).Value = 309  # I107: 8677 -> 309
$ws.Cells.Item(107, 10).Value = 1250  # J107: 100 -> 1250
$ws.Cells.Item(107, 11).Value = 309  # K107: 8677 -> 309
$ws.Cells.Item(107, 12).Value = 1250  # L107: 100 -> 1250
$ws.Cells.Item(107, 13).Value = 1611  # M107: -6757 -> 1611
$ws.Cells.Item(107, 14).Value = -5090  # N107: -3940 -> -5090

# Row 109: A Time for Peace | Smilodonskin Codex
$ws.Cells.Item(109, 8).Value = 33927  # H109: 37929.668 -> 33927
$ws.Cells.Item(109, 10).Value = 33927  # J109: 37929.668 -> 33927
$ws.Cells.Item(109, 12).Value = 33927  # L109: 37929.668 -> 33927
$ws.Cells.Item(109, 14).Value = -36701  # N109: -40703.668 -> -36701

# Row 117: A Greater Grimoire | Zonureskin Grimoire
$ws.Cells.Item(117, 8).Value = 48716.8  # H117: 48718.4 -> 48716.8
$ws.Cells.Item(117, 10).Value = 48716.8  # J117: 48718.4 -> 48716.8
$ws.Cells.Item(117, 12).Value = 48716.8  # L117: 48718.4 -> 48716.8
$ws.Cells.Item(117, 14).Value = -57894.8  # N117: -57896.4 -> -57894.8

# Row 133: Big Brush, Big Dreams | Ginseng Angle Brush
$ws.Cells.Item(133, 8).Value = 51944.223  # H133: 52383.332 -> 51944.223
$ws.Cells.Item(133, 10).Value = 51944.223  # J133: 52383.332 -> 51944.223
$ws.Cells.Item(133, 12).Value = 51944.223  # L133: 52383.332 -> 51944.223
$ws.Cells.Item(133, 14).Value = -62064.223  # N133: -62503.332 -> -62064.223

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 1692.3511  # H138: 2190.8071 -> 1692.3511
$ws.Cells.Item(138, 9).Value = 1016.86664  # I138: 0 -> 1016.86664
$ws.Cells.Item(138, 10).Value = 2312.6938  # J138: 2190.8071 -> 2312.6938
$ws.Cells.Item(138, 11).Value = 3050.59992  # K138: 0 -> 3050.59992
$ws.Cells.Item(138, 12).Value = 6938.0814  # L138: 6572.4213 -> 6938.0814
$ws.Cells.Item(138, 13).Value = 2089.40008  # M138: None -> 2089.40008
$ws.Cells.Item(138, 14).Value = -17218.0814  # N138: -16852.4213 -> -17218.0814

$ws = $wb.Worksheets.Item("ARM")
# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Cells.Item(110, 8).Value = 2298.353  # H110: 2317.8 -> 2298.353
$ws.Cells.Item(110, 9).Value = 2312.2856  # I110: 2319.0715 -> 2312.2856
$ws.Cells.Item(110, 10).Value = 2233.3333  # J110: 2300 -> 2233.3333
$ws.Cells.Item(110, 11).Value = 2312.2856  # K110: 2319.0715 -> 2312.2856
$ws.Cells.Item(110, 12).Value = 2233.3333  # L110: 2300 -> 2233.3333
$ws.Cells.Item(110, 13).Value = -267.2856000000002  # M110: -274.0715 -> -267.2856000000002
$ws.Cells.Item(110, 14).Value = -6323.3333  # N110: -6390 -> -6323.3333

# Row 114: A New Regular | Bluespirit Gauntlets of Fending
$ws.Cells.Item(114, 8).Value = 45961.332  # H114: 45964 -> 45961.332
$ws.Cells.Item(114, 10).Value = 45961.332  # J114: 45964 -> 45961.332
$ws.Cells.Item(114, 12).Value = 45961.332  # L114: 45964 -> 45961.332
$ws.Cells.Item(114, 14).Value = -54639.332  # N114: -54642 -> -54639.332

# Row 117: Signed, Shield, Delivered | Titanbronze Tower Shield
$ws.Cells.Item(117, 8).Value = 47997  # H117: 49119.5 -> 47997
$ws.Cells.Item(117, 10).Value = 47997  # J117: 49119.5 -> 47997
$ws.Cells.Item(117, 12).Value = 47997  # L117: 49119.5 -> 47997
$ws.Cells.Item(117, 14).Value = -57175  # N117: -58297.5 -> -57175

# Row 118: A Budding Business | Titanbronze Headband of Scouting
$ws.Cells.Item(118, 8).Value = 49401  # H118: 49409 -> 49401
$ws.Cells.Item(118, 10).Value = 49401  # J118: 49409 -> 49401
$ws.Cells.Item(118, 12).Value = 49401  # L118: 49409 -> 49401
$ws.Cells.Item(118, 14).Value = -52715  # N118: -52723 -> -52715

# Row 121: Shield to Shield | Dwarven Mythril Shield
$ws.Cells.Item(121, 8).Value = 34426.4  # H121: 31459 -> 34426.4
$ws.Cells.Item(121, 10).Value = 34426.4  # J121: 31459 -> 34426.4
$ws.Cells.Item(121, 12).Value = 34426.4  # L121: 31459 -> 34426.4
$ws.Cells.Item(121, 14).Value = -37920.4  # N121: -34953 -> -37920.4

# Row 131: Additions to the Armoire | Chondrite Top of Maiming
$ws.Cells.Item(131, 8).Value = 47037.332  # H131: 47042.668 -> 47037.332
$ws.Cells.Item(131, 10).Value = 47037.332  # J131: 47042.668 -> 47037.332
$ws.Cells.Item(131, 12).Value = 47037.332  # L131: 47042.668 -> 47037.332
$ws.Cells.Item(131, 14).Value = -57117.332  # N131: -57122.668 -> -57117.332

# Row 137: Odd Instruments | Cobalt Tungsten Alembic
$ws.Cells.Item(137, 8).Value = 37250  # H137: 37800 -> 37250
$ws.Cells.Item(137, 10).Value = 37250  # J137: 37800 -> 37250
$ws.Cells.Item(137, 12).Value = 37250  # L137: 37800 -> 37250
$ws.Cells.Item(137, 14).Value = -47450  # N137: -48000 -> -47450

$ws = $wb.Worksheets.Item("BSM")
# Row 57: No Refunds, Only Exchanges | Cobalt File
$ws.Cells.Item(57, 8).Value = 55419.5  # H57: 55444.5 -> 55419.5
$ws.Cells.Item(57, 10).Value = 55419.5  # J57: 55444.5 -> 55419.5
$ws.Cells.Item(57, 12).Value = 55419.5  # L57: 55444.5 -> 55419.5
$ws.Cells.Item(57, 14).Value = -56859.5  # N57: -56884.5 -> -56859.5

# Row 122: To Delight a Dancer | High Durium Tathlums
$ws.Cells.Item(122, 8).Value = 40505.75  # H122: 40469.6 -> 40505.75
$ws.Cells.Item(122, 10).Value = 40505.75  # J122: 40469.6 -> 40505.75
$ws.Cells.Item(122, 12).Value = 40505.75  # L122: 40469.6 -> 40505.75
$ws.Cells.Item(122, 14).Value = -50305.75  # N122: -50269.6 -> -50305.75

# Row 130: Annals of the Empire I | Chondrite Magitek Axe
$ws.Cells.Item(130, 8).Value = 49181  # H130: 49183.668 -> 49181
$ws.Cells.Item(130, 10).Value = 49181  # J130: 49183.668 -> 49181
$ws.Cells.Item(130, 12).Value = 49181  # L130: 49183.668 -> 49181
$ws.Cells.Item(130, 14).Value = -59221  # N130: -59223.668 -> -59221

# Row 132: Always Be Prepaired | Mountain Chromite Twinfangs
$ws.Cells.Item(132, 8).Value = 25063.5  # H132: 25335.264 -> 25063.5
$ws.Cells.Item(132, 10).Value = 25063.5  # J132: 25335.264 -> 25063.5
$ws.Cells.Item(132, 12).Value = 25063.5  # L132: 25335.264 -> 25063.5
$ws.Cells.Item(132, 14).Value = -35183.5  # N132: -35455.264 -> -35183.5

# Row 136: Maintaining the Maintainers | Cobalt Tungsten File
$ws.Cells.Item(136, 8).Value = 55419.5  # H136: 55444.5 -> 55419.5
$ws.Cells.Item(136, 10).Value = 55419.5  # J136: 55444.5 -> 55419.5
$ws.Cells.Item(136, 12).Value = 55419.5  # L136: 55444.5 -> 55419.5
$ws.Cells.Item(136, 14).Value = -65619.5  # N136: -65644.5 -> -65619.5

# Row 139: Maul Me | Titanium Gold Maul
$ws.Cells.Item(139, 8).Value = 44855.8  # H139: 45715.8 -> 44855.8
$ws.Cells.Item(139, 10).Value = 44855.8  # J139: 45715.8 -> 44855.8
$ws.Cells.Item(139, 12).Value = 44855.8  # L139: 45715.8 -> 44855.8
$ws.Cells.Item(139, 14).Value = -55135.8  # N139: -55995.8 -> -55135.8

$ws = $wb.Worksheets.Item("CRP")
# Row 43: The Long Lance of the Law | Steel Halberd
$ws.Cells.Item(43, 8).Value = 48632  # H43: 48636 -> 48632
$ws.Cells.Item(43, 10).Value = 48632  # J43: 48636 -> 48632
$ws.Cells.Item(43, 12).Value = 48632  # L43: 48636 -> 48632
$ws.Cells.Item(43, 14).Value = -49000  # N43: -49004 -> -49000

# Row 52: Spin It Like You Mean It | Mahogany Spinning Wheel
$ws.Cells.Item(52, 8).Value = 71596.664  # H52: 71630 -> 71596.664
$ws.Cells.Item(52, 10).Value = 71596.664  # J52: 71630 -> 71596.664
$ws.Cells.Item(52, 12).Value = 71596.664  # L52: 71630 -> 71596.664
$ws.Cells.Item(52, 14).Value = -72184.664  # N52: -72218 -> -72184.664

# Row 96: Composition | Larch Composite Bow
$ws.Cells.Item(96, 8).Value = 71844.8  # H96: 71850.39999999999 -> 71844.8
$ws.Cells.Item(96, 10).Value = 71844.8  # J96: 71850.39999999999 -> 71844.8
$ws.Cells.Item(96, 12).Value = 71844.8  # L96: 71850.39999999999 -> 71844.8
$ws.Cells.Item(96, 14).Value = -77336.8  # N96: -77342.39999999999 -> -77336.8

# Row 100: Run Before They Walk | Pine Cane
$ws.Cells.Item(100, 8).Value = 46972  # H100: 43436 -> 46972
$ws.Cells.Item(100, 10).Value = 46972  # J100: 43436 -> 46972
$ws.Cells.Item(100, 12).Value = 46972  # L100: 43436 -> 46972
$ws.Cells.Item(100, 14).Value = -49136  # N100: -45600 -> -49136

# Row 101: Everybody's Heard about the 'Berd | Doman Steel Halberd
$ws.Cells.Item(101, 8).Value = 48632  # H101: 48636 -> 48632
$ws.Cells.Item(101, 10).Value = 48632  # J101: 48636 -> 48632
$ws.Cells.Item(101, 12).Value = 48632  # L101: 48636 -> 48632
$ws.Cells.Item(101, 14).Value = -55122  # N101: -55126 -> -55122

# Row 111: Taking Aim | Applewood Longbow
$ws.Cells.Item(111, 8).Value = 48694  # H111: 48702 -> 48694
$ws.Cells.Item(111, 10).Value = 48694  # J111: 48702 -> 48694
$ws.Cells.Item(111, 12).Value = 48694  # L111: 48702 -> 48694
$ws.Cells.Item(111, 14).Value = -56874  # N111: -56882 -> -56874

# Row 116: The Right Tool for the Job | Sandteak Rod
$ws.Cells.Item(116, 8).Value = 49819.668  # H116: 49822.332 -> 49819.668
$ws.Cells.Item(116, 10).Value = 49819.668  # J116: 49822.332 -> 49819.668
$ws.Cells.Item(116, 12).Value = 49819.668  # L116: 49822.332 -> 49819.668
$ws.Cells.Item(116, 14).Value = -58997.668  # N116: -59000.332 -> -58997.668

# Row 131: An Integral Reward | Integral Necklace of Crafting
$ws.Cells.Item(131, 8).Value = 37496  # H131: 39964 -> 37496
$ws.Cells.Item(131, 10).Value = 37496  # J131: 39964 -> 37496
$ws.Cells.Item(131, 12).Value = 37496  # L131: 39964 -> 37496
$ws.Cells.Item(131, 14).Value = -47576  # N131: -50044 -> -47576

# Row 133: Yimepi's Country Charms | Ginseng Earrings
$ws.Cells.Item(133, 8).Value = 27997.715  # H133: 27998.857 -> 27997.715
$ws.Cells.Item(133, 10).Value = 27997.715  # J133: 27998.857 -> 27997.715
$ws.Cells.Item(133, 12).Value = 27997.715  # L133: 27998.857 -> 27997.715
$ws.Cells.Item(133, 14).Value = -33057.715  # N133: -33058.857 -> -33057.715

# Row 137: Lament of the Lazylump | Dark Mahogany Fishing Rod
$ws.Cells.Item(137, 8).Value = 38242.785  # H137: 38271.355 -> 38242.785
$ws.Cells.Item(137, 10).Value = 38242.785  # J137: 38271.355 -> 38242.785
$ws.Cells.Item(137, 12).Value = 38242.785  # L137: 38271.355 -> 38242.785
$ws.Cells.Item(137, 14).Value = -48442.785  # N137: -48471.355 -> -48442.785

# Row 139: Weaving a Path | Acacia Spinning Wheel
$ws.Cells.Item(139, 8).Value = 60439.8  # H139: 60459.8 -> 60439.8
$ws.Cells.Item(139, 10).Value = 65049.75  # J139: 65074.75 -> 65049.75
$ws.Cells.Item(139, 12).Value = 65049.75  # L139: 65074.75 -> 65049.75
$ws.Cells.Item(139, 14).Value = -75329.75  # N139: -75354.75 -> -75329.75

$ws = $wb.Worksheets.Item("GSM")
# Row 53: North Ore South | Electrum Gorget
$ws.Cells.Item(53, 8).Value = 27988.25  # H53: 24134.334 -> 27988.25
$ws.Cells.Item(53, 10).Value = 27988.25  # J53: 24134.334 -> 27988.25
$ws.Cells.Item(53, 12).Value = 27988.25  # L53: 24134.334 -> 27988.25
$ws.Cells.Item(53, 14).Value = -29250.25  # N53: -25396.334 -> -29250.25

# Row 105: Untucked | Palladium Tuck
$ws.Cells.Item(105, 8).Value = 43444.668  # H105: 46171 -> 43444.668
$ws.Cells.Item(105, 10).Value = 43444.668  # J105: 46171 -> 43444.668
$ws.Cells.Item(105, 12).Value = 43444.668  # L105: 46171 -> 43444.668
$ws.Cells.Item(105, 14).Value = -50432.668  # N105: -53159 -> -50432.668

# Row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Cells.Item(107, 8).Value = 3588.818  # H107: 2771.0667 -> 3588.818
$ws.Cells.Item(107, 9).Value = 384.625  # I107: 314.4 -> 384.625
$ws.Cells.Item(107, 10).Value = 12133.333  # J107: 3999.4 -> 12133.333
$ws.Cells.Item(107, 11).Value = 384.625  # K107: 314.4 -> 384.625
$ws.Cells.Item(107, 12).Value = 12133.333  # L107: 3999.4 -> 12133.333
$ws.Cells.Item(107, 13).Value = 1535.375  # M107: 1605.6 -> 1535.375
$ws.Cells.Item(107, 14).Value = -15973.333  # N107: -7839.4 -> -15973.333

# Row 110: Slimming Down | Stonegold Rapier
$ws.Cells.Item(110, 8).Value = 49999  # H110: 49233.332 -> 49999
$ws.Cells.Item(110, 10).Value = 49999  # J110: 49233.332 -> 49999
$ws.Cells.Item(110, 12).Value = 49999  # L110: 49233.332 -> 49999
$ws.Cells.Item(110, 14).Value = -58179  # N110: -57413.332 -> -58179

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Cells.Item(113, 8).Value = 1287.3077  # H113: 1281.6 -> 1287.3077
$ws.Cells.Item(113, 9).Value = 1158  # I113: 1238.875 -> 1158
$ws.Cells.Item(113, 10).Value = 1578.25  # J113: 1330.4286 -> 1578.25
$ws.Cells.Item(113, 11).Value = 1158  # K113: 1238.875 -> 1158
$ws.Cells.Item(113, 12).Value = 1578.25  # L113: 1330.4286 -> 1578.25
$ws.Cells.Item(113, 13).Value = 1012  # M113: 931.125 -> 1012
$ws.Cells.Item(113, 14).Value = -5918.25  # N113: -5670.4286 -> -5918.25

# Row 119: Bulking Up | Dwarven Mythril Rapier
$ws.Cells.Item(119, 8).Value = 48753  # H119: 48761 -> 48753
$ws.Cells.Item(119, 10).Value = 48753  # J119: 48761 -> 48753
$ws.Cells.Item(119, 12).Value = 48753  # L119: 48761 -> 48753
$ws.Cells.Item(119, 14).Value = -58429  # N119: -58437 -> -58429

# Row 124: The Sage's Successor | Pewter Pendulums
$ws.Cells.Item(124, 8).Value = 37587.332  # H124: 37920.668 -> 37587.332
$ws.Cells.Item(124, 10).Value = 37587.332  # J124: 37920.668 -> 37587.332
$ws.Cells.Item(124, 12).Value = 37587.332  # L124: 37920.668 -> 37587.332
$ws.Cells.Item(124, 14).Value = -47407.332  # N124: -47740.668 -> -47407.332

# Row 137: Sew Excited | Cobalt Tungsten Needle
$ws.Cells.Item(137, 8).Value = 41855.8  # H137: 42415.8 -> 41855.8
$ws.Cells.Item(137, 10).Value = 41855.8  # J137: 42415.8 -> 41855.8
$ws.Cells.Item(137, 12).Value = 41855.8  # L137: 42415.8 -> 41855.8
$ws.Cells.Item(137, 14).Value = -52055.8  # N137: -52615.8 -> -52055.8

# Row 139: Ringing Gratitude | White Gold Ring of Healing
$ws.Cells.Item(139, 8).Value = 33972  # H139: 50144 -> 33972
$ws.Cells.Item(139, 10).Value = 33972  # J139: 50144 -> 33972
$ws.Cells.Item(139, 12).Value = 33972  # L139: 50144 -> 33972
$ws.Cells.Item(139, 14).Value = -44252  # N139: -60424 -> -44252

$ws = $wb.Worksheets.Item("LTW")
# Row 36: Campaign in the Membrane | Toadskin Jacket
$ws.Cells.Item(36, 8).Value = 48663.25  # H36: 48667.75 -> 48663.25
$ws.Cells.Item(36, 10).Value = 48663.25  # J36: 48667.75 -> 48663.25
$ws.Cells.Item(36, 12).Value = 48663.25  # L36: 48667.75 -> 48663.25
$ws.Cells.Item(36, 14).Value = -49787.25  # N36: -49791.75 -> -49787.25

# Row 108: Girding for Glory | Smilodonskin Trousers of Maiming
$ws.Cells.Item(108, 8).Value = 0  # H108: 48618 -> 0
$ws.Cells.Item(108, 10).Value = 0  # J108: 48618 -> 0
$ws.Cells.Item(108, 12).Value = 0  # L108: 48618 -> 0
$ws.Cells.Item(108, 14).ClearContents()  # remove N108 (was -56298)

# Row 110: Breeches of Trust | Gliderskin Breeches of Fending
$ws.Cells.Item(110, 8).Value = 38363.75  # H110: 38518.75 -> 38363.75
$ws.Cells.Item(110, 10).Value = 38363.75  # J110: 38518.75 -> 38363.75
$ws.Cells.Item(110, 12).Value = 38363.75  # L110: 38518.75 -> 38363.75
$ws.Cells.Item(110, 14).Value = -46543.75  # N110: -46698.75 -> -46543.75

# Row 112: A Slippery Slope | Gliderskin Boots of Casting
$ws.Cells.Item(112, 8).Value = 43157.332  # H112: 43164 -> 43157.332
$ws.Cells.Item(112, 10).Value = 43157.332  # J112: 43164 -> 43157.332
$ws.Cells.Item(112, 12).Value = 43157.332  # L112: 43164 -> 43157.332
$ws.Cells.Item(112, 14).Value = -46111.332  # N112: -46118 -> -46111.332

# Row 119: Fit for a Friend | Swallowskin Gloves of Fending
$ws.Cells.Item(119, 8).Value = 44896  # H119: 45530.668 -> 44896
$ws.Cells.Item(119, 10).Value = 44896  # J119: 45530.668 -> 44896
$ws.Cells.Item(119, 12).Value = 44896  # L119: 45530.668 -> 44896
$ws.Cells.Item(119, 14).Value = -54572  # N119: -55206.668 -> -54572

# Row 131: For What Was Gleaned | Ophiotauroskin Wristband of Gathering
$ws.Cells.Item(131, 8).Value = 45296.668  # H131: 46875 -> 45296.668
$ws.Cells.Item(131, 10).Value = 45296.668  # J131: 46875 -> 45296.668
$ws.Cells.Item(131, 12).Value = 45296.668  # L131: 46875 -> 45296.668
$ws.Cells.Item(131, 14).Value = -55376.668  # N131: -56955 -> -55376.668

# Row 133: The Perfect Accessory | Loboskin Amulet of Fending
$ws.Cells.Item(133, 8).Value = 39833.332  # H133: 38618 -> 39833.332
$ws.Cells.Item(133, 10).Value = 39833.332  # J133: 38618 -> 39833.332
$ws.Cells.Item(133, 12).Value = 39833.332  # L133: 38618 -> 39833.332
$ws.Cells.Item(133, 14).Value = -44893.332  # N133: -43678 -> -44893.332

# Row 137: Lending Artisans a Hand | Br'aaxskin Halfgloves of Crafting
$ws.Cells.Item(137, 8).Value = 41350  # H137: 41416.668 -> 41350
$ws.Cells.Item(137, 10).Value = 41350  # J137: 41416.668 -> 41350
$ws.Cells.Item(137, 12).Value = 41350  # L137: 41416.668 -> 41350
$ws.Cells.Item(137, 14).Value = -51550  # N137: -51616.668 -> -51550

$ws = $wb.Worksheets.Item("WVR")
# Row 119: A Job Well Done | Dwarven Cotton Gaskins of Fending
$ws.Cells.Item(119, 8).Value = 49690  # H119: 49698 -> 49690
$ws.Cells.Item(119, 10).Value = 49690  # J119: 49698 -> 49690
$ws.Cells.Item(119, 12).Value = 49690  # L119: 49698 -> 49690
$ws.Cells.Item(119, 14).Value = -59366  # N119: -59374 -> -59366

# Row 127: Turban Sprawl | Snow Linen Turban of Crafting
$ws.Cells.Item(127, 8).Value = 42421  # H127: 42429 -> 42421
$ws.Cells.Item(127, 10).Value = 42421  # J127: 42429 -> 42421
$ws.Cells.Item(127, 12).Value = 42421  # L127: 42429 -> 42421
$ws.Cells.Item(127, 14).Value = -52341  # N127: -52349 -> -52341

# Row 128: Lightening Up | Scarlet Moko Gaskins of the Rising Dragon
$ws.Cells.Item(128, 8).Value = 48569  # H128: 42676.75 -> 48569
$ws.Cells.Item(128, 10).Value = 48569  # J128: 42676.75 -> 48569
$ws.Cells.Item(128, 12).Value = 48569  # L128: 42676.75 -> 48569
$ws.Cells.Item(128, 14).Value = -58529  # N128: -52636.75 -> -58529

# Row 139: Cruel Climates | Rroneek Serge Trousers of Gathering
$ws.Cells.Item(139, 8).Value = 52199.8  # H139: 52619.8 -> 52199.8
$ws.Cells.Item(139, 10).Value = 52199.8  # J139: 52619.8 -> 52199.8
$ws.Cells.Item(139, 12).Value = 52199.8  # L139: 52619.8 -> 52199.8
$ws.Cells.Item(139, 14).Value = -62479.8  # N139: -62899.8 -> -62479.8
